# Update "want to go" counts (column F) across all four sheets
# (展览 / 演出 / 本地生活 / 全部类型) to match the refreshed data pull
# recorded in commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# 展览
$ws1.Range("F6").Value = 2654
$ws1.Range("F10").Value = 11
$ws1.Range("F11").Value = 1536
$ws1.Range("F14").Value = 638
$ws1.Range("F16").Value = 1385
$ws1.Range("F17").Value = 22
$ws1.Range("F18").Value = 9
$ws1.Range("F19").Value = 546
$ws1.Range("F20").Value = 3943
$ws1.Range("F21").Value = 3943
$ws1.Range("F22").Value = 656
$ws1.Range("F23").Value = 3306
$ws1.Range("F24").Value = 763
$ws1.Range("F25").Value = 25
$ws1.Range("F26").Value = 2204
$ws1.Range("F28").Value = 308
$ws1.Range("F30").Value = 29
$ws1.Range("F31").Value = 1168
$ws1.Range("F32").Value = 765
$ws1.Range("F34").Value = 1050
$ws1.Range("F35").Value = 1050

# 演出
$ws2.Range("F12").Value = 103
$ws2.Range("F18").Value = 256
$ws2.Range("F19").Value = 196

# 本地生活
$ws3.Range("F4").Value = 538
$ws3.Range("F5").Value = 145

# 全部类型
$ws4.Range("F9").Value = 538
$ws4.Range("F10").Value = 2654
$ws4.Range("F11").Value = 2654
$ws4.Range("F20").Value = 11
$ws4.Range("F23").Value = 1536
$ws4.Range("F27").Value = 103
$ws4.Range("F28").Value = 1385
$ws4.Range("F29").Value = 22
$ws4.Range("F30").Value = 546
$ws4.Range("F32").Value = 3943
$ws4.Range("F33").Value = 3943
$ws4.Range("F34").Value = 656
$ws4.Range("F35").Value = 3306
$ws4.Range("F36").Value = 763
$ws4.Range("F37").Value = 2204
$ws4.Range("F39").Value = 308
$ws4.Range("F41").Value = 29
$ws4.Range("F42").Value = 1168
$ws4.Range("F44").Value = 256
$ws4.Range("F45").Value = 196
$ws4.Range("F47").Value = 765
$ws4.Range("F49").Value = 1050
$ws4.Range("F50").Value = 1050
